$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A21").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A35").Value = "lang_bus_transport"
$ws.Range("B35").Value = "Xe đưa đón"
$ws.Range("C35").Value = "Bus No"

$ws.Range("C35").Select()
